$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "28.660.29"
$ws.Range("E2").Value = "  -3.20%  "
$ws.Range("D3").Value = "1.851.67"
$ws.Range("E4").Value = "  -1.13%  "
$ws.Range("D5").Value = "'334.61"
$ws.Range("E5").Value = "  +2.45%  "
$ws.Range("E6").Value = "  -0.94%  "
$ws.Range("D7").Value = "'0.4652"
$ws.Range("E7").Value = "  -3.63%  "
$ws.Range("D8").Value = "'0.3909"
$ws.Range("E8").Value = "  -3.83%  "
$ws.Range("D9").Value = "'46.18"
$ws.Range("E9").Value = "  -3.14%  "
$ws.Range("D10").Value = "'0.07915"
$ws.Range("E10").Value = "  -3.96%  "
$ws.Range("D11").Value = "'0.9846"
$ws.Range("E11").Value = "  -2.76%  "
$ws.Range("D12").Value = "'22.31"
$ws.Range("E12").Value = "  -6.54%  "
$ws.Range("D13").Value = "1.837.44"
$ws.Range("E13").Value = "  -4.79%  "
$ws.Range("D14").Value = "'5.851"
$ws.Range("E14").Value = "  -4.33%  "
$ws.Range("D15").Value = "'7.009"
$ws.Range("E15").Value = "  -3.72%  "
$ws.Range("D16").Value = "'0.06855"
$ws.Range("E16").Value = "  -0.14%  "
$ws.Range("D17").Value = "'87.75"
$ws.Range("E17").Value = "  -4.56%  "
$ws.Range("E18").Value = "  -1.01%  "
$ws.Range("D19").Value = "'0.00001008"
$ws.Range("E19").Value = "  -3.09%  "
$ws.Range("D20").Value = "'17.14"
$ws.Range("E20").Value = "  -3.04%  "
$ws.Range("E21").Value = "  -0.84%  "
$ws.Range("D22").Value = "28.675.85"
$ws.Range("E22").Value = "  -3.13%  "
$ws.Range("D23").Value = "'5.401"
$ws.Range("E23").Value = "  -5.17%  "
$ws.Range("E24").Value = "  -5.88%  "
$ws.Range("D25").Value = "'2.133"
$ws.Range("E25").Value = "  -2.56%  "
$ws.Range("D26").Value = "2.089.07"
$ws.Range("E26").Value = "  -3.34%  "
$ws.Range("D27").Value = "'153.17"
$ws.Range("E27").Value = "  -2.00%  "
$ws.Range("D28").Value = "'19.49"
$ws.Range("E28").Value = "  -2.89%  "
$ws.Range("D29").Value = "'6.053"
$ws.Range("E29").Value = "  -6.05%  "
$ws.Range("E30").Value = "  -3.15%  "
$ws.Range("D31").Value = "'117.42"
$ws.Range("E31").Value = "  -2.86%  "
$ws.Range("D32").Value = "'0.9792"
$ws.Range("E32").Value = "  -3.68%  "
$ws.Range("D33").Value = "'0.09422"
$ws.Range("E33").Value = "  -2.36%  "
$ws.Range("D34").Value = "'5.380"
$ws.Range("E34").Value = "  -4.49%  "
$ws.Range("D35").Value = "'3.484"
$ws.Range("E35").Value = "  -2.59%  "
$ws.Range("D36").Value = "'1.351"
$ws.Range("E36").Value = "  -2.14%  "
$ws.Range("D37").Value = "'0.06158"
$ws.Range("E37").Value = "  -3.71%  "
$ws.Range("E38").Value = "  -4.20%  "
$ws.Range("D39").Value = "'1.173"
$ws.Range("E39").Value = "  -1.29%  "
$ws.Range("D40").Value = "'1.001"
$ws.Range("E40").Value = "  -0.98%  "
$ws.Range("D41").Value = "'0.5726"
$ws.Range("E41").Value = "  -4.06%  "
$ws.Range("D42").Value = "'7.615"
$ws.Range("E42").Value = "  -3.67%  "
$ws.Range("E43").Value = "  -5.08%  "
$ws.Range("E44").Value = "  -2.84%  "
$ws.Range("D45").Value = "'2.372"
$ws.Range("E45").Value = "  -2.43%  "
$ws.Range("D46").Value = "'1.245"
$ws.Range("E46").Value = "  -3.09%  "
$ws.Range("D47").Value = "'0.5393"
$ws.Range("E47").Value = "  -3.28%  "
$ws.Range("D48").Value = "'11.82"
$ws.Range("E48").Value = "  -5.22%  "
$ws.Range("D49").Value = "'0.07136"
$ws.Range("E49").Value = "  -5.49%  "
$ws.Range("E50").Value = "  -3.26%  "
$ws.Range("D51").Value = "'115.47"
$ws.Range("E51").Value = "  -3.48%  "
